$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in C1 (renamed from "NOME+EXTENSÃO*" to "NOME_ANEXO+EXTENSÃO*")
$ws.Range("C1").Value = "NOME_ANEXO+EXTENSÃO*"

# Give the header an underline in addition to the existing bold/white styling
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").Font.Underline = $true
$ws.Range("C1").Font.ThemeColor = 0

# Leave C1 selected, matching the saved view state
$ws.Range("C1").Select()
